# Auto-generated edit script: updates Leve profit-calculation columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets to refreshed market-price data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 11.666667
$ws.Range("I9").Value = 11.666667
$ws.Range("K9").Value = 11.666667
$ws.Range("M9").Value = 157.333333
$ws.Range("H17").Value = 688.6667
$ws.Range("J17").Value = 688.6667
$ws.Range("L17").Value = 2066.0001
$ws.Range("N17").Value = -2402.0001
$ws.Range("H40").Value = 3497
$ws.Range("J40").Value = 3621.25
$ws.Range("L40").Value = 3621.25
$ws.Range("N40").Value = -3971.25
$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H49").Value = 100
$ws.Range("I49").Value = 100
$ws.Range("K49").Value = 300
$ws.Range("M49").Value = -164
$ws.Range("H51").Value = 9331
$ws.Range("J51").Value = 9499.5
$ws.Range("L51").Value = 9499.5
$ws.Range("N51").Value = -10467.5
$ws.Range("H58").Value = 2529.625
$ws.Range("I58").Value = 1164
$ws.Range("J58").Value = 4805.6665
$ws.Range("K58").Value = 3492
$ws.Range("L58").Value = 14416.9995
$ws.Range("M58").Value = -3342
$ws.Range("N58").Value = -14716.9995
$ws.Range("H62").Value = 4875
$ws.Range("I62").Value = 4750
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4750
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4126
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4875
$ws.Range("I65").Value = 4750
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 23750
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -20630
$ws.Range("N65").Value = -31240
$ws.Range("H106").Value = 25005000
$ws.Range("I106").Value = 25005000
$ws.Range("K106").Value = 25005000
$ws.Range("M106").Value = -25004369
$ws.Range("H137").Value = 2299.75
$ws.Range("J137").Value = 2399
$ws.Range("L137").Value = 7197
$ws.Range("N137").Value = -12297
$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5099.4
$ws.Range("I21").Value = 4999
$ws.Range("J21").Value = 5124.5
$ws.Range("K21").Value = 4999
$ws.Range("L21").Value = 5124.5
$ws.Range("M21").Value = -4625
$ws.Range("N21").Value = -5872.5
$ws.Range("H23").Value = 252007
$ws.Range("J23").Value = 252007
$ws.Range("L23").Value = 252007
$ws.Range("N23").Value = -252525
$ws.Range("H32").Value = 1136.6
$ws.Range("I32").Value = 766.5833
$ws.Range("K32").Value = 766.5833
$ws.Range("M32").Value = -479.5833
$ws.Range("H61").Value = 840.25
$ws.Range("I61").Value = 840.25
$ws.Range("K61").Value = 840.25
$ws.Range("M61").Value = -628.25
$ws.Range("H63").Value = 4269.8
$ws.Range("I63").Value = 2837.25
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 2837.25
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -2151.25
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 4269.8
$ws.Range("I66").Value = 2837.25
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 14186.25
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -10754.25
$ws.Range("N66").Value = -56864
$ws.Range("H74").Value = 1472.6666
$ws.Range("J74").Value = 1168
$ws.Range("L74").Value = 1168
$ws.Range("N74").Value = -2916
$ws.Range("H77").Value = 1472.6666
$ws.Range("J77").Value = 1168
$ws.Range("L77").Value = 5840
$ws.Range("N77").Value = -14576
$ws.Range("H122").Value = 1995
$ws.Range("I122").Value = 492.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 1477.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = 972.5
$ws.Range("N122").Value = -19900
$ws.Range("H136").Value = 840.25
$ws.Range("I136").Value = 840.25
$ws.Range("K136").Value = 2520.75
$ws.Range("M136").Value = 29.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H86").Value = 2413
$ws.Range("I86").Value = 2666.25
$ws.Range("K86").Value = 2666.25
$ws.Range("M86").Value = -1543.25
$ws.Range("H89").Value = 2413
$ws.Range("I89").Value = 2666.25
$ws.Range("K89").Value = 13331.25
$ws.Range("M89").Value = -7715.25
$ws.Range("H105").Value = 8999.666999999999
$ws.Range("I105").Value = 8999.666999999999
$ws.Range("K105").Value = 8999.666999999999
$ws.Range("M105").Value = -7252.666999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -1700
$ws.Range("H37").Value = 10000
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10214
$ws.Range("H58").Value = 11475
$ws.Range("I58").Value = 6750
$ws.Range("J58").Value = 12000
$ws.Range("K58").Value = 6750
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -6547
$ws.Range("N58").Value = -12406
$ws.Range("H134").Value = 8000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 11475
$ws.Range("I136").Value = 6750
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 20250
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -17700
$ws.Range("N136").Value = -41100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 986.0909
$ws.Range("I34").Value = 86.85714
$ws.Range("J34").Value = 2559.75
$ws.Range("K34").Value = 260.57142
$ws.Range("L34").Value = 7679.25
$ws.Range("M34").Value = -176.57142
$ws.Range("N34").Value = -7847.25
$ws.Range("H39").Value = 4516.8887
$ws.Range("J39").Value = 5099.8
$ws.Range("L39").Value = 15299.4
$ws.Range("N39").Value = -15887.4
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H55").Value = 4685.7144
$ws.Range("I55").Value = 800
$ws.Range("J55").Value = 6240
$ws.Range("K55").Value = 2400
$ws.Range("L55").Value = 18720
$ws.Range("M55").Value = -2223
$ws.Range("N55").Value = -19074
$ws.Range("H69").Value = 1808.6
$ws.Range("J69").Value = 1808.6
$ws.Range("L69").Value = 5425.799999999999
$ws.Range("N69").Value = -7047.799999999999
$ws.Range("H72").Value = 1808.6
$ws.Range("J72").Value = 1808.6
$ws.Range("L72").Value = 16277.4
$ws.Range("N72").Value = -24389.4
$ws.Range("H131").Value = 1979.1538
$ws.Range("J131").Value = 1979.1538
$ws.Range("L131").Value = 5937.4614
$ws.Range("N131").Value = -16017.4614
$ws.Range("H140").Value = 1866.6666
$ws.Range("I140").Value = 1800
$ws.Range("J140").Value = 1900
$ws.Range("K140").Value = 5400
$ws.Range("L140").Value = 5700
$ws.Range("M140").Value = -220
$ws.Range("N140").Value = -16060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 10051.5
$ws.Range("I22").Value = 10051.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 10051.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -9522.5
$ws.Range("N22").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H43").Value = 3984.5
$ws.Range("I43").Value = 3984.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3984.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3833.5
$ws.Range("N43").ClearContents()
$ws.Range("H122").Value = 3056.3333
$ws.Range("I122").Value = 3316
$ws.Range("J122").Value = 2796.6667
$ws.Range("K122").Value = 9948
$ws.Range("L122").Value = 8390.000100000001
$ws.Range("M122").Value = -7498
$ws.Range("N122").Value = -13290.0001
$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -104899
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999.5
$ws.Range("I7").Value = 1999.5
$ws.Range("K7").Value = 1999.5
$ws.Range("M7").Value = -1887.5
$ws.Range("H46").Value = 1750
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -2376
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 1999.5
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5
